$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, applied to both "展览" and "全部类型" sheets
$updates = @{
    2  = 1167
    3  = 124
    4  = 1637
    5  = 631
    8  = 11587
    11 = 456
    12 = 373
    14 = 813
    15 = 12412
    16 = 13157
    18 = 145
    21 = 246
    24 = 133
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
